$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the style (bold font, thin box border, center/top alignment) on B1
$cellB1 = $ws.Range("B1")
$cellB1.Font.Bold = $true
$cellB1.HorizontalAlignment = -4108
$cellB1.VerticalAlignment = -4160
$cellB1.Borders.LineStyle = 1
$cellB1.Borders.Weight = 2

# Copy the same formatting onto A2 so both cells reuse the same style record
$cellB1.Copy()
$cellA2 = $ws.Range("A2")
$cellA2.PasteSpecial(-4122)

$excel.CutCopyMode = $false
